$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "abg_hypercap_threshold", "Symptom – Respiratory", 2384, 10749, 22.18, 21.39, 22.96)
    ,@(3, "abg_hypercap_threshold", "Symptom – Nervous", 1595, 10749, 14.84, 14.17, 15.51)
    ,@(4, "abg_hypercap_threshold", "Symptom – Digestive", 1259, 10749, 11.71, 11.1, 12.32)
    ,@(5, "abg_hypercap_threshold", "Injuries & adverse effects", 1031, 10749, 9.59, 9.03, 10.15)
    ,@(6, "abg_hypercap_threshold", "Symptom – Circulatory", 982, 10749, 9.14, 8.59, 9.68)
    ,@(7, "abg_hypercap_threshold", "Diseases (patient-stated)", 950, 10749, 8.84, 8.3, 9.37)
    ,@(8, "abg_hypercap_threshold", "Other", 944, 10749, 8.78, 8.25, 9.32)
    ,@(9, "abg_hypercap_threshold", "Symptom – Genitourinary", 627, 10749, 5.83, 5.39, 6.28)
    ,@(10, "abg_hypercap_threshold", "Symptom – General", 520, 10749, 4.84, 4.43, 5.24)
    ,@(11, "abg_hypercap_threshold", "Symptom – Skin/Hair/Nails", 285, 10749, 2.65, 2.35, 2.96)
    ,@(12, "abg_hypercap_threshold", "Symptom – Eye/Ear", 172, 10749, 1.6, 1.36, 1.84)
    ,@(13, "any_hypercap_icd", "Symptom – Respiratory", 997, 1983, 50.28, 48.08, 52.48)
    ,@(14, "any_hypercap_icd", "Symptom – Nervous", 232, 1983, 11.7, 10.28, 13.11)
    ,@(15, "any_hypercap_icd", "Symptom – Circulatory", 137, 1983, 6.91, 5.79, 8.02)
    ,@(16, "any_hypercap_icd", "Symptom – Digestive", 133, 1983, 6.71, 5.61, 7.81)
    ,@(17, "any_hypercap_icd", "Other", 100, 1983, 5.04, 4.08, 6.01)
    ,@(18, "any_hypercap_icd", "Injuries & adverse effects", 86, 1983, 4.34, 3.44, 5.23)
    ,@(19, "any_hypercap_icd", "Symptom – General", 83, 1983, 4.19, 3.3, 5.07)
    ,@(20, "any_hypercap_icd", "Symptom – Genitourinary", 82, 1983, 4.14, 3.26, 5.01)
    ,@(21, "any_hypercap_icd", "Diseases (patient-stated)", 62, 1983, 3.13, 2.36, 3.89)
    ,@(22, "any_hypercap_icd", "Symptom – Skin/Hair/Nails", 37, 1983, 1.87, 1.27, 2.46)
    ,@(23, "any_hypercap_icd", "Symptom – Eye/Ear", 34, 1983, 1.71, 1.14, 2.29)
    ,@(24, "other_hypercap_threshold", "Symptom – Respiratory", 6879, 36518, 18.84, 18.44, 19.24)
    ,@(25, "other_hypercap_threshold", "Symptom – Nervous", 5896, 36518, 16.15, 15.77, 16.52)
    ,@(26, "other_hypercap_threshold", "Symptom – Digestive", 4926, 36518, 13.49, 13.14, 13.84)
    ,@(27, "other_hypercap_threshold", "Symptom – Circulatory", 3774, 36518, 10.33, 10.02, 10.65)
    ,@(28, "other_hypercap_threshold", "Diseases (patient-stated)", 3069, 36518, 8.4, 8.12, 8.69)
    ,@(29, "other_hypercap_threshold", "Other", 2933, 36518, 8.03, 7.75, 8.31)
    ,@(30, "other_hypercap_threshold", "Injuries & adverse effects", 2760, 36518, 7.56, 7.29, 7.83)
    ,@(31, "other_hypercap_threshold", "Symptom – Genitourinary", 2289, 36518, 6.27, 6.02, 6.52)
    ,@(32, "other_hypercap_threshold", "Symptom – General", 2265, 36518, 6.2, 5.96, 6.45)
    ,@(33, "other_hypercap_threshold", "Symptom – Skin/Hair/Nails", 879, 36518, 2.41, 2.25, 2.56)
    ,@(34, "other_hypercap_threshold", "Symptom – Eye/Ear", 848, 36518, 2.32, 2.17, 2.48)
    ,@(35, "pco2_threshold_any", "Symptom – Respiratory", 7995, 41080, 19.46, 19.08, 19.84)
    ,@(36, "pco2_threshold_any", "Symptom – Nervous", 6562, 41080, 15.97, 15.62, 16.33)
    ,@(37, "pco2_threshold_any", "Symptom – Digestive", 5765, 41080, 14.03, 13.7, 14.37)
    ,@(38, "pco2_threshold_any", "Symptom – Circulatory", 4159, 41080, 10.12, 9.83, 10.42)
    ,@(39, "pco2_threshold_any", "Diseases (patient-stated)", 3315, 41080, 8.07, 7.81, 8.33)
    ,@(40, "pco2_threshold_any", "Other", 3208, 41080, 7.81, 7.55, 8.07)
    ,@(41, "pco2_threshold_any", "Injuries & adverse effects", 2920, 41080, 7.11, 6.86, 7.36)
    ,@(42, "pco2_threshold_any", "Symptom – Genitourinary", 2586, 41080, 6.3, 6.06, 6.53)
    ,@(43, "pco2_threshold_any", "Symptom – General", 2570, 41080, 6.26, 6.02, 6.49)
    ,@(44, "pco2_threshold_any", "Symptom – Eye/Ear", 1022, 41080, 2.49, 2.34, 2.64)
    ,@(45, "pco2_threshold_any", "Symptom – Skin/Hair/Nails", 978, 41080, 2.38, 2.23, 2.53)
    ,@(46, "vbg_hypercap_threshold", "Symptom – Respiratory", 4457, 17548, 25.4, 24.75, 26.04)
    ,@(47, "vbg_hypercap_threshold", "Symptom – Digestive", 2612, 17548, 14.88, 14.36, 15.41)
    ,@(48, "vbg_hypercap_threshold", "Symptom – Nervous", 2459, 17548, 14.01, 13.5, 14.53)
    ,@(49, "vbg_hypercap_threshold", "Symptom – Circulatory", 1689, 17548, 9.63, 9.19, 10.06)
    ,@(50, "vbg_hypercap_threshold", "Other", 1130, 17548, 6.44, 6.08, 6.8)
    ,@(51, "vbg_hypercap_threshold", "Symptom – General", 1106, 17548, 6.3, 5.94, 6.66)
    ,@(52, "vbg_hypercap_threshold", "Injuries & adverse effects", 1100, 17548, 6.27, 5.91, 6.63)
    ,@(53, "vbg_hypercap_threshold", "Diseases (patient-stated)", 1065, 17548, 6.07, 5.72, 6.42)
    ,@(54, "vbg_hypercap_threshold", "Symptom – Genitourinary", 994, 17548, 5.66, 5.32, 6.01)
    ,@(55, "vbg_hypercap_threshold", "Symptom – Eye/Ear", 556, 17548, 3.17, 2.91, 3.43)
    ,@(56, "vbg_hypercap_threshold", "Symptom – Skin/Hair/Nails", 380, 17548, 2.17, 1.95, 2.38)
)

foreach ($r in $data) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
}